# Update the practice-sheet division problems in the single table.
# Each content row (1, 5, 9, 13, 17) holds 5 problems (columns 1-5);
# the intervening rows are blank spacer rows.
#
# Cells are updated by assigning Range.Text directly (rather than
# Find/Replace) because Find.Execute on a Cell's Range in this runtime
# is not scoped to that cell and can match/replace other identical
# occurrences elsewhere in the document.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "34÷9="
$t.Cell(1, 2).Range.Text = "85÷6="
$t.Cell(1, 3).Range.Text = "39÷6="
$t.Cell(1, 4).Range.Text = "40÷4="
$t.Cell(1, 5).Range.Text = "33÷9="
$t.Cell(5, 1).Range.Text = "98÷7="
$t.Cell(5, 2).Range.Text = "67÷4="
$t.Cell(5, 3).Range.Text = "16÷4="
$t.Cell(5, 4).Range.Text = "88÷9="
$t.Cell(5, 5).Range.Text = "84÷7="
$t.Cell(9, 1).Range.Text = "65÷3="
$t.Cell(9, 2).Range.Text = "37÷4="
$t.Cell(9, 3).Range.Text = "15÷9="
$t.Cell(9, 4).Range.Text = "73÷2="
$t.Cell(9, 5).Range.Text = "53÷9="
$t.Cell(13, 1).Range.Text = "95÷9="
$t.Cell(13, 2).Range.Text = "42÷5="
$t.Cell(13, 3).Range.Text = "32÷6="
$t.Cell(13, 4).Range.Text = "50÷3="
$t.Cell(13, 5).Range.Text = "25÷2="
$t.Cell(17, 1).Range.Text = "10÷7="
$t.Cell(17, 2).Range.Text = "66÷8="
$t.Cell(17, 3).Range.Text = "10÷8="
$t.Cell(17, 4).Range.Text = "81÷9="
$t.Cell(17, 5).Range.Text = "17÷8="
